$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the dataset. It belongs right
# after the existing row for this market/category ordering, so insert a
# fresh row at position 420 (this shifts the former rows 420:491 down to
# 421:492, preserving all of their data untouched).
$ws.Rows(420).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(420, 1).Value = 6
$ws.Cells.Item(420, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(420, 3).Value = "Metropolitana"
$ws.Cells.Item(420, 4).Value = 44474
$ws.Cells.Item(420, 5).Value = 13
$ws.Cells.Item(420, 6).Value = 100112031
$ws.Cells.Item(420, 7).Value = "Poroto verde"
$ws.Cells.Item(420, 8).Value = "Magnum"
$ws.Cells.Item(420, 9).Value = "Primera"
$ws.Cells.Item(420, 10).Value = 250
$ws.Cells.Item(420, 11).Value = 40000
$ws.Cells.Item(420, 12).Value = 45000
$ws.Cells.Item(420, 13).Value = 43000
$ws.Cells.Item(420, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(420, 15).Value = "Perú"
$ws.Cells.Item(420, 16).Value = 1720
$ws.Cells.Item(420, 17).Value = 25
$ws.Cells.Item(420, 18).Value = "Hortaliza"
